$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42605.888032407405
$ws.Range("B3").Value = -30
$ws.Range("C3").Value = 52
$ws.Range("D3").Value = 46
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 86
$ws.Range("G3").Value = 23420
$ws.Range("H3").Value = 12671
$ws.Range("I3").Value = 1383
$ws.Range("J3").Value = 146
$ws.Range("K3").Value = 128
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 32
$ws.Range("N3").Value = "Bag"
